# Edit script implementing the commit: "revised text to refer to help"
$p = $ppt.ActivePresentation

# --- Slide 11 ---
$s11 = $p.Slides.Item(11)

# Shape 2 "TextBox 4": extend the suggestion bullet to reference the Help page.
$sh11_2 = $s11.Shapes.Item(2)
$tr11_2 = $sh11_2.TextFrame.TextRange
$para9 = $tr11_2.Paragraphs(9,1)

$full9 = $para9.Characters(1, $para9.Length)
$full9.Text = " change message to: ‘’No scientific name was found. Please try again.  The Help page has tips on getting good photos.’’"

$r1 = $para9.Characters(1, 68)
$r1.Text = " change message to: ‘’No scientific name was found. Please try again"

$r2 = $para9.Characters(69, 21)
$r2.Text = ".  The Help page has "

$r3 = $para9.Characters(90, 28)
$r3.Text = "tips on getting good photos."

$r4 = $para9.Characters(118, 1)
$r4.Text = "’"

$r5 = $para9.Characters(119, 1)
$r5.Text = "’"

# Shape 3 "TextBox 6": widen the box (to fit the longer title) and reword the title.
$sh11_3 = $s11.Shapes.Item(3)
$sh11_3.Left = 242.598
$sh11_3.Width = 510.8563

$tr11_3 = $sh11_3.TextFrame.TextRange
$title = $tr11_3.Paragraphs(1,1)
$titleFull = $title.Characters(1, $title.Length)
$titleFull.Text = "9. Instruct user in empty capture alert"

# --- Slide 15 ---
$s15 = $p.Slides.Item(15)

# Shape 1 picture: move down to make room for the extra title line below.
$sh15_1 = $s15.Shapes.Item(1)
$sh15_1.Top = 104.1354

# Shape 15 "TextBox 61": reposition/resize and add a heading line above "Priority = low".
$sh15_15 = $s15.Shapes.Item(15)
$sh15_15.Left = 50.6372
$sh15_15.Top = 21.0218
$sh15_15.Width = 311.56954
$sh15_15.Height = 84.82032

$tr15_15 = $sh15_15.TextFrame.TextRange
[void]$tr15_15.InsertBefore("Alert for no-text-found`r")

$para2 = $tr15_15.Paragraphs(2,1)
$pr1 = $para2.Characters(1,9)
$pr1.Text = "Priority "
$pr2 = $para2.Characters(10,2)
$pr2.Text = "= "
$pr3 = $para2.Characters(12,3)
$pr3.Text = "low"
